$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44: the debit amount was corrected (larger order placed on 2-Feb-2020)
$ws.Range("B44").Value = 53040
# D44 formatting was nudged to match the "Ordered Amount" style used nearby
$ws.Range("D44").Style = $ws.Range("D42").Style

# Row 45: newly recorded transaction on 2-Feb-2020 (credit / "Manual Added")
$ws.Range("A45").Value = 43863
$ws.Range("C45").Value = 105590
$ws.Range("D45").Formula = "=D43"
$ws.Range("D45").Style = $ws.Range("D43").Style

# Row 46: newly recorded transaction on 3-Feb-2020 (debit / "Ordered Amount")
$ws.Range("A46").Value = 43864
$ws.Range("B46").Value = 11440
$ws.Range("D46").Formula = "=D44"
$ws.Range("D46").Style = $ws.Range("D44").Style

# Keep the frozen-pane view scrolled to the newly entered rows
$ws.Range("E48").Select()
